$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D = Price, Column E = Volume(1h)
# Refresh cryptos list figures (scheduled GitHub Actions data pull).
# Values are written with a temporary text ("@") number format so that
# digit-and-dot strings such as "292.01" or "6.040" are kept as literal
# text (matching the source feed) instead of being auto-converted to
# numbers by Excel; the style is then reset back to Normal/General so
# the cell formatting matches the original (unstyled) cells.

# Row 2
$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "22.476.60"
$c.Style = "Normal"
$c = $ws.Cells.Item(2, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.34%  "
$c.Style = "Normal"

# Row 3
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "1.572.41"
$c.Style = "Normal"
$c = $ws.Cells.Item(3, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.09%  "
$c.Style = "Normal"

# Row 4
$c = $ws.Cells.Item(4, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.01%  "
$c.Style = "Normal"

# Row 5
$c = $ws.Cells.Item(5, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.05%  "
$c.Style = "Normal"

# Row 6
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "292.01"
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.16%  "
$c.Style = "Normal"

# Row 7
$c = $ws.Cells.Item(7, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.07%  "
$c.Style = "Normal"

# Row 8
$c = $ws.Cells.Item(8, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.28%  "
$c.Style = "Normal"

# Row 9
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.3399"
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.49%  "
$c.Style = "Normal"

# Row 10
$c = $ws.Cells.Item(10, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.35%  "
$c.Style = "Normal"

# Row 11
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.07545"
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.98%  "
$c.Style = "Normal"

# Row 12
$c = $ws.Cells.Item(12, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.02%  "
$c.Style = "Normal"

# Row 13
$c = $ws.Cells.Item(13, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.32%  "
$c.Style = "Normal"

# Row 14
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "6.040"
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.82%  "
$c.Style = "Normal"

# Row 15
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "6.957"
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.33%  "
$c.Style = "Normal"

# Row 16
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "1.574.31"
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.54%  "
$c.Style = "Normal"

# Row 17
$c = $ws.Cells.Item(17, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.89%  "
$c.Style = "Normal"

# Row 18
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "90.77"
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.58%  "
$c.Style = "Normal"

# Row 19
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "0.06761"
$c.Style = "Normal"
$c = $ws.Cells.Item(19, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.28%  "
$c.Style = "Normal"

# Row 20
$c = $ws.Cells.Item(20, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.02%  "
$c.Style = "Normal"

# Row 21
$c = $ws.Cells.Item(21, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.72%  "
$c.Style = "Normal"

# Row 22
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "16.35"
$c.Style = "Normal"
$c = $ws.Cells.Item(22, 5)
$c.NumberFormat = "@"
$c.Value = "  -2.38%  "
$c.Style = "Normal"

# Row 23
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "12.18"
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.61%  "
$c.Style = "Normal"

# Row 24
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "22.477.48"
$c.Style = "Normal"
$c = $ws.Cells.Item(24, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.37%  "
$c.Style = "Normal"

# Row 25
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "2.380"
$c.Style = "Normal"
$c = $ws.Cells.Item(25, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.57%  "
$c.Style = "Normal"

# Row 26
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "2.618"
$c.Style = "Normal"
$c = $ws.Cells.Item(26, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.80%  "
$c.Style = "Normal"

# Row 27
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "20.02"
$c.Style = "Normal"
$c = $ws.Cells.Item(27, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.45%  "
$c.Style = "Normal"

# Row 28
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "149.16"
$c.Style = "Normal"
$c = $ws.Cells.Item(28, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.36%  "
$c.Style = "Normal"

# Row 29
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "5.039"
$c.Style = "Normal"
$c = $ws.Cells.Item(29, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.17%  "
$c.Style = "Normal"

# Row 30
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "125.32"
$c.Style = "Normal"
$c = $ws.Cells.Item(30, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.97%  "
$c.Style = "Normal"

# Row 31
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "1.746.77"
$c.Style = "Normal"
$c = $ws.Cells.Item(31, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.23%  "
$c.Style = "Normal"

# Row 32
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "1.084"
$c.Style = "Normal"
$c = $ws.Cells.Item(32, 5)
$c.NumberFormat = "@"
$c.Value = "  +9.11%  "
$c.Style = "Normal"

# Row 33
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "6.201"
$c.Style = "Normal"
$c = $ws.Cells.Item(33, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.71%  "
$c.Style = "Normal"

# Row 34
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "2.013"
$c.Style = "Normal"
$c = $ws.Cells.Item(34, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.08%  "
$c.Style = "Normal"

# Row 35
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "9.806"
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 5)
$c.NumberFormat = "@"
$c.Value = "  -3.29%  "
$c.Style = "Normal"

# Row 36
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "0.08372"
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.56%  "
$c.Style = "Normal"

# Row 37
$c = $ws.Cells.Item(37, 5)
$c.NumberFormat = "@"
$c.Value = "  -2.02%  "
$c.Style = "Normal"

# Row 38
$c = $ws.Cells.Item(38, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.23%  "
$c.Style = "Normal"

# Row 39
$c = $ws.Cells.Item(39, 5)
$c.NumberFormat = "@"
$c.Value = "  -2.28%  "
$c.Style = "Normal"

# Row 40
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "0.06536"
$c.Style = "Normal"
$c = $ws.Cells.Item(40, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.62%  "
$c.Style = "Normal"

# Row 41
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "5.445"
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.78%  "
$c.Style = "Normal"

# Row 42
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "11.34"
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.40%  "
$c.Style = "Normal"

# Row 43
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "0.6231"
$c.Style = "Normal"
$c = $ws.Cells.Item(43, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.78%  "
$c.Style = "Normal"

# Row 45
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "14.07"
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.67%  "
$c.Style = "Normal"

# Row 46
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "3.815"
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.63%  "
$c.Style = "Normal"

# Row 47
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "0.5843"
$c.Style = "Normal"
$c = $ws.Cells.Item(47, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.83%  "
$c.Style = "Normal"

# Row 48
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "130.57"
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 5)
$c.NumberFormat = "@"
$c.Value = "  +4.61%  "
$c.Style = "Normal"

# Row 49
$c = $ws.Cells.Item(49, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.96%  "
$c.Style = "Normal"

# Row 50
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "1.213"
$c.Style = "Normal"
$c = $ws.Cells.Item(50, 5)
$c.NumberFormat = "@"
$c.Value = "  -4.51%  "
$c.Style = "Normal"

# Row 51
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "0.07326"
$c.Style = "Normal"
$c = $ws.Cells.Item(51, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.05%  "
$c.Style = "Normal"
